$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 88151.5
$ws.Range("I69").Value = 65006
$ws.Range("J69").Value = 99724.25
$ws.Range("K69").Value = 195018
$ws.Range("L69").Value = 299172.75
$ws.Range("M69").Value = -194144
$ws.Range("N69").Value = -300920.75
$ws.Range("H70").Value = 7315.8
$ws.Range("J70").Value = 7524.0835
$ws.Range("L70").Value = 22572.2505
$ws.Range("N70").Value = -23112.2505
$ws.Range("H72").Value = 88151.5
$ws.Range("I72").Value = 65006
$ws.Range("J72").Value = 99724.25
$ws.Range("K72").Value = 585054
$ws.Range("L72").Value = 897518.25
$ws.Range("M72").Value = -580686
$ws.Range("N72").Value = -906254.25
$ws.Range("H73").Value = 7315.8
$ws.Range("J73").Value = 7524.0835
$ws.Range("L73").Value = 22572.2505
$ws.Range("N73").Value = -24444.2505
$ws.Range("H80").Value = 3573.361
$ws.Range("I80").Value = 2207.7144
$ws.Range("J80").Value = 4442.409
$ws.Range("K80").Value = 6623.1432
$ws.Range("L80").Value = 13327.227
$ws.Range("M80").Value = -5625.1432
$ws.Range("N80").Value = -15323.227
$ws.Range("H83").Value = 3573.361
$ws.Range("I83").Value = 2207.7144
$ws.Range("J83").Value = 4442.409
$ws.Range("K83").Value = 19869.4296
$ws.Range("L83").Value = 39981.681
$ws.Range("M83").Value = -14877.4296
$ws.Range("N83").Value = -49965.681
$ws.Range("H112").Value = 1897.4286
$ws.Range("J112").Value = 1897.6154
$ws.Range("L112").Value = 5692.8462
$ws.Range("N112").Value = -7908.8462
$ws.Range("H132").Value = 12362.8545
$ws.Range("I132").Value = 998.87177
$ws.Range("J132").Value = 61606.777
$ws.Range("K132").Value = 2996.61531
$ws.Range("L132").Value = 184820.331
$ws.Range("M132").Value = -466.6153100000001
$ws.Range("N132").Value = -189880.331
$ws.Range("H137").Value = 8345597
$ws.Range("I137").Value = 14304324
$ws.Range("J137").Value = 3379.6
$ws.Range("K137").Value = 42912972
$ws.Range("L137").Value = 10138.8
$ws.Range("M137").Value = -42910422
$ws.Range("N137").Value = -15238.8
$ws.Range("H138").Value = 5132.4814
$ws.Range("I138").Value = 799.5
$ws.Range("K138").Value = 2398.5
$ws.Range("M138").Value = 2741.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 190
$ws.Range("I4").Value = 190
$ws.Range("K4").Value = 190
$ws.Range("M4").Value = -74
$ws.Range("H32").Value = 6903.712
$ws.Range("J32").Value = 2525.5
$ws.Range("L32").Value = 2525.5
$ws.Range("N32").Value = -3099.5
$ws.Range("H74").Value = 1730.8
$ws.Range("I74").Value = 1620.1765
$ws.Range("J74").Value = 1965.875
$ws.Range("K74").Value = 1620.1765
$ws.Range("L74").Value = 1965.875
$ws.Range("M74").Value = -746.1765
$ws.Range("N74").Value = -3713.875
$ws.Range("H77").Value = 1730.8
$ws.Range("I77").Value = 1620.1765
$ws.Range("J77").Value = 1965.875
$ws.Range("K77").Value = 8100.8825
$ws.Range("L77").Value = 9829.375
$ws.Range("M77").Value = -3732.8825
$ws.Range("N77").Value = -18565.375
$ws.Range("H122").Value = 3949.1667
$ws.Range("I122").Value = 3586.7585
$ws.Range("J122").Value = 4369.56
$ws.Range("K122").Value = 10760.2755
$ws.Range("L122").Value = 13108.68
$ws.Range("M122").Value = -8310.2755
$ws.Range("N122").Value = -18008.68
$ws.Range("H132").Value = 5262.151
$ws.Range("I132").Value = 4581.2744
$ws.Range("K132").Value = 13743.8232
$ws.Range("M132").Value = -11213.8232

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 95000
$ws.Range("J59").Value = 95000
$ws.Range("L59").Value = 95000
$ws.Range("N59").Value = -96694
$ws.Range("H87").Value = 30909
$ws.Range("J87").Value = 31500
$ws.Range("L87").Value = 31500
$ws.Range("N87").Value = -33996
$ws.Range("H90").Value = 30909
$ws.Range("J90").Value = 31500
$ws.Range("L90").Value = 94500
$ws.Range("N90").Value = -106980
$ws.Range("H94").Value = 6871.9165
$ws.Range("I94").Value = 5769.125
$ws.Range("J94").Value = 9077.5
$ws.Range("K94").Value = 5769.125
$ws.Range("L94").Value = 9077.5
$ws.Range("M94").Value = -5318.125
$ws.Range("N94").Value = -9979.5
$ws.Range("H134").Value = 5190.6665
$ws.Range("I134").Value = 3038.5908
$ws.Range("K134").Value = 9115.7724
$ws.Range("M134").Value = -6580.7724

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4047.0386
$ws.Range("I31").Value = 4360.75
$ws.Range("J31").Value = 3341.1875
$ws.Range("K31").Value = 4360.75
$ws.Range("L31").Value = 3341.1875
$ws.Range("M31").Value = -4065.75
$ws.Range("N31").Value = -3931.1875
$ws.Range("H34").Value = 4047.0386
$ws.Range("I34").Value = 4360.75
$ws.Range("J34").Value = 3341.1875
$ws.Range("K34").Value = 4360.75
$ws.Range("L34").Value = 3341.1875
$ws.Range("M34").Value = -4158.75
$ws.Range("N34").Value = -3745.1875
$ws.Range("H35").Value = 1794.3
$ws.Range("I35").Value = 880.5
$ws.Range("J35").Value = 5449.5
$ws.Range("K35").Value = 880.5
$ws.Range("L35").Value = 5449.5
$ws.Range("M35").Value = -586.5
$ws.Range("N35").Value = -6037.5
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()
$ws.Range("H99").Value = 37400704
$ws.Range("I99").Value = 8132739.5
$ws.Range("K99").Value = 8132739.5
$ws.Range("M99").Value = -8131241.5
$ws.Range("H105").Value = 60187690
$ws.Range("I105").Value = 37040036
$ws.Range("K105").Value = 37040036
$ws.Range("M105").Value = -37038289
$ws.Range("H122").Value = 932114.9399999999
$ws.Range("I122").Value = 1138018.2
$ws.Range("J122").Value = 5550
$ws.Range("K122").Value = 3414054.6
$ws.Range("L122").Value = 16650
$ws.Range("M122").Value = -3411604.6
$ws.Range("N122").Value = -21550
$ws.Range("H126").Value = 37400704
$ws.Range("I126").Value = 8132739.5
$ws.Range("K126").Value = 24398218.5
$ws.Range("M126").Value = -24395748.5
$ws.Range("H134").Value = 2554.6287
$ws.Range("I134").Value = 2585.9707
$ws.Range("K134").Value = 7757.9121
$ws.Range("M134").Value = -5222.9121

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 170.61539
$ws.Range("I6").Value = 170.61539
$ws.Range("K6").Value = 511.84617
$ws.Range("M6").Value = -398.84617
$ws.Range("H59").Value = 34853.668
$ws.Range("J59").Value = 41704.9
$ws.Range("L59").Value = 125114.7
$ws.Range("N59").Value = -126194.7
$ws.Range("H68").Value = 1383.3334
$ws.Range("J68").Value = 1500
$ws.Range("L68").Value = 4500
$ws.Range("N68").Value = -6122
$ws.Range("H71").Value = 1383.3334
$ws.Range("J71").Value = 1500
$ws.Range("L71").Value = 13500
$ws.Range("N71").Value = -21612

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2859.3684
$ws.Range("I102").Value = 2445.7942
$ws.Range("K102").Value = 2445.7942
$ws.Range("M102").Value = -823.7941999999998
$ws.Range("H107").Value = 393.75
$ws.Range("I107").Value = 375
$ws.Range("K107").Value = 375
$ws.Range("M107").Value = 1545
$ws.Range("H113").Value = 9658.883
$ws.Range("I113").Value = 8333.166999999999
$ws.Range("J113").Value = 10382
$ws.Range("K113").Value = 8333.166999999999
$ws.Range("L113").Value = 10382
$ws.Range("M113").Value = -6163.166999999999
$ws.Range("N113").Value = -14722
$ws.Range("H122").Value = 5533.9565
$ws.Range("I122").Value = 4821.2
$ws.Range("J122").Value = 6082.231
$ws.Range("K122").Value = 14463.6
$ws.Range("L122").Value = 18246.693
$ws.Range("M122").Value = -12013.6
$ws.Range("N122").Value = -23146.693
$ws.Range("H132").Value = 6753.22
$ws.Range("I132").Value = 6140.521
$ws.Range("K132").Value = 18421.563
$ws.Range("M132").Value = -15891.563

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 871.4286
$ws.Range("J22").Value = 999.5
$ws.Range("L22").Value = 999.5
$ws.Range("N22").Value = -1589.5
$ws.Range("H27").Value = 871.4286
$ws.Range("J27").Value = 999.5
$ws.Range("L27").Value = 999.5
$ws.Range("N27").Value = -1213.5
$ws.Range("H40").Value = 2797.4167
$ws.Range("I40").Value = 2400.4443
$ws.Range("K40").Value = 2400.4443
$ws.Range("M40").Value = -2264.4443
$ws.Range("H46").Value = 3906.1035
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 3906.1035
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 3906.1035
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -4282.1035
$ws.Range("H132").Value = 4539.96
$ws.Range("I132").Value = 4563
$ws.Range("J132").Value = 4275
$ws.Range("K132").Value = 13689
$ws.Range("L132").Value = 12825
$ws.Range("M132").Value = -11159
$ws.Range("N132").Value = -17885
$ws.Range("H136").Value = 3467.2
$ws.Range("I136").Value = 3230.6924
$ws.Range("K136").Value = 9692.0772
$ws.Range("M136").Value = -7142.0772

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1472.25
$ws.Range("I122").Value = 1472.25
$ws.Range("K122").Value = 4416.75
$ws.Range("M122").Value = -1966.75
$ws.Range("H132").Value = 2529.9473
$ws.Range("I132").Value = 2740
$ws.Range("K132").Value = 8220
$ws.Range("M132").Value = -5690
$ws.Range("H133").Value = 99635
$ws.Range("J133").Value = 99635
$ws.Range("L133").Value = 99635
$ws.Range("N133").Value = -109755
$ws.Range("H136").Value = 4041.56
$ws.Range("I136").Value = 3628.9119
$ws.Range("J136").Value = 4918.4375
$ws.Range("K136").Value = 10886.7357
$ws.Range("L136").Value = 14755.3125
$ws.Range("M136").Value = -8336.735700000001
$ws.Range("N136").Value = -19855.3125
